$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Create Item Category -> Update Unit)
$ws.Name = "Update Unit"

# Update the header text and drop the now-unused second column
$ws.Range("A1").Value = "s"
$ws.Range("B:B").EntireColumn.Delete()
